$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B-column (Actual Production) values for rows 2-29 (January production data)
$newB = @{
    2  = 1643
    3  = 1561
    4  = 1487
    5  = 1402
    6  = 1377
    7  = 1365
    8  = 1328
    9  = 1314
    10 = 1380
    11 = 1397
    12 = 1396
    13 = 1341
    14 = 1307
    15 = 1305
    16 = 1272
    17 = 1175
    18 = 1060
    19 = 998
    20 = 942
    21 = 815
    22 = 678
    23 = 568
    24 = 503
    25 = 456
    26 = 425
    27 = 384
    28 = 364
    29 = 349
}

# Shift every timestamp in column A (rows 2-97) forward by 2 days,
# and update column B values for rows 2-29 with the new production figures.
for ($r = 2; $r -le 97; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $d = $cellA.Value()
    $cellA.Value = $d.AddDays(2)

    if ($newB.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $newB[$r]
    }
}
